# "New tenant support in live" - append new sprint-run rows to the
# AMSIN, BETA and AMS history sheets (and backfill formatting that was
# missing on the previous last row of AMSIN, row 75).

$wb = $excel.ActiveWorkbook

function Write-HistoryRow {
    param(
        $ws,
        [int]$row,
        [string]$dateText,
        [double]$timeSerial,
        [string]$sprintName,
        [double]$total,
        [double]$pass,
        [double]$fail,
        [double]$timeTaken,
        [int]$formatSourceRow
    )

    # Column A ("Run Date") is stored as literal text even though it
    # looks like a date (e.g. "2023-04-18"). Clear any previous content
    # first (so a pre-existing cell picks up fresh formatting the same
    # way a brand new cell would) then flip to Text number-format
    # before writing, so Excel's date auto-detection doesn't turn the
    # string into a date serial.
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Clear()
    $aCell.NumberFormat = "@"
    $aCell.Value = $dateText

    # Column B ("Run Time") keeps the custom date-time display format
    # used by the rest of the column - copy it from an existing,
    # already-formatted cell so we don't invent a new number format.
    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Clear()
    $ws.Cells.Item($formatSourceRow, 2).Copy()
    $bCell.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $bCell.Value = $timeSerial

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Clear()
    $cCell.Value = $sprintName

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Clear()
    $dCell.Value = $total

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Clear()
    $eCell.Value = $pass

    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Clear()
    $fCell.Value = $fail

    $gCell = $ws.Cells.Item($row, 7)
    $gCell.Clear()
    $gCell.Value = $timeTaken
}

# ---------------------------------------------------------------------
# AMSIN: backfill row 75 (the previous final row, which was missing the
# usual formatting) and append the two new rows, 76 and 77.
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Write-HistoryRow $wsAmsin 75 "2023-04-18" 45034.60912424768 "176firsttrail" 165 165 0 4.93 74
Write-HistoryRow $wsAmsin 76 "2023-04-19" 45035.71066576389 "176scndcyc"    165 164 1 5.09 74
Write-HistoryRow $wsAmsin 77 "2023-04-20" 45036.42039621528 "176fnlruntest" 165 165 0 5.16 74

# ---------------------------------------------------------------------
# BETA: append row 35.
# ---------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

Write-HistoryRow $wsBeta 35 "2023-04-20" 45036.52886015047 "176beta" 165 165 0 5.13 34

# ---------------------------------------------------------------------
# AMS: append rows 41 and 42.
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

Write-HistoryRow $wsAms 41 "2023-05-04" 45050.72406798611 "176firstsycle" 165 165 0 4.47 40
Write-HistoryRow $wsAms 42 "2023-05-08" 45054.54967030385 "176htfxtrl"    165 165 0 4.58 40
